$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to remain plain text so numeric-looking
# strings (e.g. "1.00", "67.264.63") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '67.264.63'
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').Value = '2.620.78'
$ws.Range('E3').Value = '  -1.14%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '591.55'
$ws.Range('E5').Value = '  -1.37%  '
$ws.Range('D6').Value = '166.71'
$ws.Range('E6').Value = '  -0.40%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '0.533'
$ws.Range('E8').Value = '  -2.14%  '
$ws.Range('D9').Value = '2.619.47'
$ws.Range('E9').Value = '  -1.15%  '
$ws.Range('E10').Value = '  -5.52%  '
$ws.Range('E11').Value = '  +1.63%  '
$ws.Range('D12').Value = '0.365'
$ws.Range('E12').Value = '  -0.63%  '
$ws.Range('E13').Value = '  -0.44%  '
$ws.Range('D14').Value = '27.37'
$ws.Range('E14').Value = '  -2.58%  '
$ws.Range('D15').Value = '3.098.93'
$ws.Range('E15').Value = '  -1.12%  '
$ws.Range('D16').Value = '0.0000180'
$ws.Range('E16').Value = '  -2.61%  '
$ws.Range('D17').Value = '67.269.11'
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('D18').Value = '2.615.43'
$ws.Range('E18').Value = '  -0.90%  '
$ws.Range('D19').Value = '11.84'
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('D20').Value = '7.89'
$ws.Range('E20').Value = '  -0.35%  '
$ws.Range('D21').Value = '355.41'
$ws.Range('E21').Value = '  -2.86%  '
$ws.Range('D22').Value = '4.30'
$ws.Range('E22').Value = '  -2.65%  '
$ws.Range('D23').Value = '4.65'
$ws.Range('E23').Value = '  -3.66%  '
$ws.Range('B24').Value = 'Aptos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D24').Value = '10.41'
$ws.Range('E24').Value = '  -2.33%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('B26').Value = 'SuiNetwork'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D26').Value = '1.93'
$ws.Range('E26').Value = '  -4.36%  '
$ws.Range('D27').Value = '69.12'
$ws.Range('E27').Value = '  -2.37%  '
$ws.Range('D28').Value = '2.755.38'
$ws.Range('E28').Value = '  -0.79%  '
$ws.Range('D29').Value = '1.01'
$ws.Range('E29').Value = '  +1.60%  '
$ws.Range('E30').Value = '  -2.58%  '
$ws.Range('D31').Value = '543.02'
$ws.Range('E31').Value = '  -2.80%  '
$ws.Range('D32').Value = '7.94'
$ws.Range('E32').Value = '  -1.57%  '
$ws.Range('E33').Value = '  -3.75%  '
$ws.Range('D34').Value = '1.88'
$ws.Range('E34').Value = '  -2.51%  '
$ws.Range('D35').Value = '0.135'
$ws.Range('E35').Value = '  +2.71%  '
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').Value = '1.50'
$ws.Range('E37').Value = '  -3.01%  '
$ws.Range('D38').Value = '157.05'
$ws.Range('E38').Value = '  -0.04%  '
$ws.Range('D39').Value = '18.93'
$ws.Range('E39').Value = '  -2.85%  '
$ws.Range('E40').Value = '  -2.23%  '
$ws.Range('E41').Value = '  +1.62%  '
$ws.Range('D42').Value = '1.81'
$ws.Range('E42').Value = '  -1.40%  '
$ws.Range('D43').Value = '5.17'
$ws.Range('E43').Value = '  -2.18%  '
$ws.Range('E45').Value = '  -4.27%  '
$ws.Range('D46').Value = '0.0₆0304'
$ws.Range('E46').Value = '  -0.14%  '
$ws.Range('D47').Value = '151.99'
$ws.Range('E47').Value = '  -1.34%  '
$ws.Range('E48').Value = '  -2.81%  '
$ws.Range('D49').Value = '3.79'
$ws.Range('E49').Value = '  -2.79%  '
$ws.Range('D50').Value = '1.71'
$ws.Range('E50').Value = '  -2.03%  '
$ws.Range('E51').Value = '  -1.36%  '
